# Update worksheet: replace "Additional Context" column with "Answer" column,
# fill in actual answers for each query, and clear the now-unused 5th column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("D1").Value = "Answer"
$ws.Range("E1").Value = $null

# Answers for each row (2-11)
$answers = @(
    "Final exams are scheduled from December 10th to December 20th.",
    "ECE 101 is an introductory course on electronics. It covers basic circuit theory and electronic devices. Instructor: Prof. Srinivasa Rao.",
    "You can access your previous semester results via the student portal under the `"Grades`" section.",
    "The library is open from 8 AM to 9 PM on weekdays, and from 10 AM to 2 PM on weekends.",
    "The prerequisite for ECE 202 is ECE 101.",
    "The instructor for CS101 is Dr. Priya Mehta.",
    "Yes, you can change your major by submitting a request to the academic office. The deadline for major change requests is November 15th.",
    "Classes for the next semester start on January 5th.",
    "You can apply for graduation by filling out the graduation application form on the student portal. The deadline for applications is April 1st.",
    "You can join a student club by attending the club fairs held at the start of each semester or by contacting the club coordinators directly."
)

$rowHeights = @(58, 87, 72.5, 58, 43.5, 43.5, 87, 43.5, 101.5, 87)

for ($i = 0; $i -lt $answers.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $answers[$i]
    $ws.Cells.Item($row, 5).Value = $null
    $ws.Rows.Item($row).RowHeight = $rowHeights[$i]
}

# Restore the view to show the top of the sheet, matching where the
# editor left the cursor after making the changes.
$ws.Range("A1").Select() | Out-Null
$ws.Range("E3").Select() | Out-Null

